$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix up the existing "2: System is very verbose..." paragraph.
#    In the source document this sentence is split across two runs by
#    a stray "_GoBack" bookmark ("...control" / " this."). Find/Replace
#    across that split merges the runs back into one and drops the
#    stale bookmark, leaving exactly the text we want.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Needs some flags to control this.", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "Needs some flags to control this.", 2) | Out-Null

# Locate that paragraph again (it is the last paragraph in the body).
$bugPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# ------------------------------------------------------------------
# 2. Add the new "2a: ..." paragraph right after it, as two runs.
# ------------------------------------------------------------------
$bugPara.Range.InsertParagraphAfter()
$p2a = $d.Paragraphs.Item($bugPara.Index + 1)
$p2aRange = $p2a.Range
$p2aRange.InsertAfter("2a: added -v flag to control verbosity.")

$p2aTail = $p2a.Range
$p2aTail.Collapse(0)
$p2aTail.InsertAfter(" Could also do with a -debug flag to control debug settings without altering program text.")

# ------------------------------------------------------------------
# 3. Add a blank paragraph after that.
# ------------------------------------------------------------------
$p2a = $d.Paragraphs.Item($p2a.Index)
$p2a.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# 4. Add the "3: debug ..." paragraph after the blank one.
# ------------------------------------------------------------------
$blankPara = $d.Paragraphs.Item($p2a.Index + 1)
$blankPara.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($blankPara.Index + 1)

$p3text = "3: debug (public Boolean) now place in AbstractSynapse and in AbstractCompartment so can be " + `
    [char]0x2019 + "d in  the Abstract constructor in  both cases. Still needs a parameter to be set up. " + `
    "Could have -debug " + [char]0x201C + "synapse" + [char]0x201D + "/" + [char]0x201D + "compartment" + `
    [char]0x201D + "/" + [char]0x201D + "all" + [char]0x201D + "  and follow this through in the constructors. Needs work."
$p3text = "3: debug (public Boolean) now place in AbstractSynapse and in AbstractCompartment so can be init" + `
    [char]0x2019 + "d in  the Abstract constructor in  both cases. Still needs a parameter to be set up. " + `
    "Could have -debug " + [char]0x201C + "synapse" + [char]0x201D + "/" + [char]0x201D + "compartment" + `
    [char]0x201D + "/" + [char]0x201D + "all" + [char]0x201D + "  and follow this through in the constructors. Needs work."

$p3.Range.InsertAfter($p3text)

# ------------------------------------------------------------------
# 5. Re-create the "_GoBack" bookmark at the very end of this last
#    paragraph (zero-width, i.e. bookmarkStart immediately followed by
#    bookmarkEnd with no text between them -- matching what Word does
#    to mark the last edit position).
#
#    A genuinely zero-width Range confuses this host's Bookmarks.Add
#    (it falls back to wrapping the whole paragraph), so instead we:
#      a) type a short placeholder at the end of the paragraph,
#      b) bookmark exactly that placeholder's (non-empty) range,
#      c) delete the placeholder text, leaving the bookmark tags
#         sitting adjacently right where the placeholder used to be.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item($p3.Index)
$tail = $p3.Range.Duplicate
$tail.Collapse(0)
$markerStart = $tail.Start
$marker = "GoBackMarker"
$tail.InsertAfter($marker)

$markerRange = $d.Range($markerStart, $markerStart + $marker.Length)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($markerStart, $markerStart + $marker.Length)
$markerRange2.Text = ""
